$wb = $excel.ActiveWorkbook

# --- Sheet 1: Neg_Change ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "LT"
$ws.Cells.Item(2, 2).Value = 3738
$ws.Cells.Item(2, 3).Value = 3740
$ws.Cells.Item(2, 4).Value = 3679.6
$ws.Cells.Item(2, 5).Value = 3704.3
$ws.Cells.Item(2, 6).Value = 2599630
$ws.Cells.Item(2, 7).Value = 5870301
$ws.Cells.Item(2, 8).Value = -0.5571555870814802
$ws.Cells.Item(2, 9).Value = "LT"
$ws.Cells.Item(3, 1).Value = "ICICIGI"
$ws.Cells.Item(3, 2).Value = 1894
$ws.Cells.Item(3, 3).Value = 1910
$ws.Cells.Item(3, 4).Value = 1871.9
$ws.Cells.Item(3, 5).Value = 1900
$ws.Cells.Item(3, 6).Value = 560303
$ws.Cells.Item(3, 7).Value = 1278304
$ws.Cells.Item(3, 8).Value = -0.5616825105765139
$ws.Cells.Item(3, 9).Value = "ICICIGI"
$ws.Cells.Item(4, 1).Value = "IDEA"
$ws.Cells.Item(4, 2).Value = 8.06
$ws.Cells.Item(4, 3).Value = 8.38
$ws.Cells.Item(4, 4).Value = 8.02
$ws.Cells.Item(4, 5).Value = 8.25
$ws.Cells.Item(4, 6).Value = 739155377
$ws.Cells.Item(4, 7).Value = 1597372619
$ws.Cells.Item(4, 8).Value = -0.5372680311355706
$ws.Cells.Item(4, 9).Value = "IDEA"
$ws.Cells.Item(5, 1).Value = "BIOCON"
$ws.Cells.Item(5, 2).Value = 340.2
$ws.Cells.Item(5, 3).Value = 343.5
$ws.Cells.Item(5, 4).Value = 338.2
$ws.Cells.Item(5, 5).Value = 341.2
$ws.Cells.Item(5, 6).Value = 1703064
$ws.Cells.Item(5, 7).Value = 3635613
$ws.Cells.Item(5, 8).Value = -0.5315607024179967
$ws.Cells.Item(5, 9).Value = "BIOCON"
$ws.Cells.Item(6, 1).Value = "LAURUSLABS"
$ws.Cells.Item(6, 2).Value = 833.4
$ws.Cells.Item(6, 3).Value = 849.8
$ws.Cells.Item(6, 4).Value = 826.2
$ws.Cells.Item(6, 5).Value = 840
$ws.Cells.Item(6, 6).Value = 2764351
$ws.Cells.Item(6, 7).Value = 5612259
$ws.Cells.Item(6, 8).Value = -0.5074441503857894
$ws.Cells.Item(6, 9).Value = "LAURUSLABS"
$ws.Cells.Item(7, 1).Value = "CROMPTON"
$ws.Cells.Item(7, 2).Value = 294.9
$ws.Cells.Item(7, 3).Value = 295.6
$ws.Cells.Item(7, 4).Value = 292.1
$ws.Cells.Item(7, 5).Value = 293.95
$ws.Cells.Item(7, 6).Value = 2414846
$ws.Cells.Item(7, 7).Value = 5080459
$ws.Cells.Item(7, 8).Value = -0.5246795614333272
$ws.Cells.Item(7, 9).Value = "CROMPTON"

# Remove the now-obsolete trailing rows 8-11 (sheet shrinks from 11 to 7 data rows)
$ws.Range("A8:I11").Delete()

# --- Sheet 2: Pos_Change ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "TECHM"
$ws.Cells.Item(2, 2).Value = 1415
$ws.Cells.Item(2, 3).Value = 1425
$ws.Cells.Item(2, 4).Value = 1391.9
$ws.Cells.Item(2, 5).Value = 1410
$ws.Cells.Item(2, 6).Value = 2708118
$ws.Cells.Item(2, 7).Value = 1851067
$ws.Cells.Item(2, 8).Value = 0.4630037702579107
$ws.Cells.Item(2, 9).Value = "TECHM"
$ws.Cells.Item(3, 1).Value = "RELIANCE"
$ws.Cells.Item(3, 2).Value = 1381.6
$ws.Cells.Item(3, 3).Value = 1389
$ws.Cells.Item(3, 4).Value = 1368
$ws.Cells.Item(3, 5).Value = 1376
$ws.Cells.Item(3, 6).Value = 14231999
$ws.Cells.Item(3, 7).Value = 9879109
$ws.Cells.Item(3, 8).Value = 0.440615646613475
$ws.Cells.Item(3, 9).Value = "RELIANCE"
$ws.Cells.Item(4, 1).Value = "COALINDIA"
$ws.Cells.Item(4, 2).Value = 389.8
$ws.Cells.Item(4, 3).Value = 391.6
$ws.Cells.Item(4, 4).Value = 385.6
$ws.Cells.Item(4, 5).Value = 388.2
$ws.Cells.Item(4, 6).Value = 6606409
$ws.Cells.Item(4, 7).Value = 4698693
$ws.Cells.Item(4, 8).Value = 0.4060099265902241
$ws.Cells.Item(4, 9).Value = "COALINDIA"
$ws.Cells.Item(5, 1).Value = "BHARTIARTL"
$ws.Cells.Item(5, 2).Value = 1928
$ws.Cells.Item(5, 3).Value = 1928
$ws.Cells.Item(5, 4).Value = 1897.5
$ws.Cells.Item(5, 5).Value = 1906.5
$ws.Cells.Item(5, 6).Value = 5072608
$ws.Cells.Item(5, 7).Value = 3239306
$ws.Cells.Item(5, 8).Value = 0.5659551768187383
$ws.Cells.Item(5, 9).Value = "BHARTIARTL"
$ws.Cells.Item(6, 1).Value = "IOC"
$ws.Cells.Item(6, 2).Value = 145.04
$ws.Cells.Item(6, 3).Value = 149.95
$ws.Cells.Item(6, 4).Value = 145.04
$ws.Cells.Item(6, 5).Value = 149.27
$ws.Cells.Item(6, 6).Value = 16918652
$ws.Cells.Item(6, 7).Value = 11925305
$ws.Cells.Item(6, 8).Value = 0.4187185988115188
$ws.Cells.Item(6, 9).Value = "IOC"
$ws.Cells.Item(7, 1).Value = "GAIL"
$ws.Cells.Item(7, 2).Value = 172.78
$ws.Cells.Item(7, 3).Value = 177.1
$ws.Cells.Item(7, 4).Value = 171.9
$ws.Cells.Item(7, 5).Value = 176.11
$ws.Cells.Item(7, 6).Value = 12460294
$ws.Cells.Item(7, 7).Value = 8760356
$ws.Cells.Item(7, 8).Value = 0.422350187595116
$ws.Cells.Item(7, 9).Value = "GAIL"
$ws.Cells.Item(8, 1).Value = "BANKBARODA"
$ws.Cells.Item(8, 2).Value = 248.5
$ws.Cells.Item(8, 3).Value = 255.39
$ws.Cells.Item(8, 4).Value = 248.31
$ws.Cells.Item(8, 5).Value = 254
$ws.Cells.Item(8, 6).Value = 11711402
$ws.Cells.Item(8, 7).Value = 7652782
$ws.Cells.Item(8, 8).Value = 0.5303456965061856
$ws.Cells.Item(8, 9).Value = "BANKBARODA"
$ws.Cells.Item(9, 1).Value = "CGPOWER"
$ws.Cells.Item(9, 2).Value = 742.35
$ws.Cells.Item(9, 3).Value = 756.4
$ws.Cells.Item(9, 4).Value = 741.75
$ws.Cells.Item(9, 5).Value = 748
$ws.Cells.Item(9, 6).Value = 2929682
$ws.Cells.Item(9, 7).Value = 1925106
$ws.Cells.Item(9, 8).Value = 0.5218289278616346
$ws.Cells.Item(9, 9).Value = "CGPOWER"
$ws.Cells.Item(10, 1).Value = "JINDALSTEL"
$ws.Cells.Item(10, 2).Value = 1035.1
$ws.Cells.Item(10, 3).Value = 1052
$ws.Cells.Item(10, 4).Value = 1035.1
$ws.Cells.Item(10, 5).Value = 1040
$ws.Cells.Item(10, 6).Value = 1278211
$ws.Cells.Item(10, 7).Value = 858507
$ws.Cells.Item(10, 8).Value = 0.4888766195266899
$ws.Cells.Item(10, 9).Value = "JINDALSTEL"
$ws.Cells.Item(11, 1).Value = "GODREJPROP"
$ws.Cells.Item(11, 2).Value = 1974
$ws.Cells.Item(11, 3).Value = 2018
$ws.Cells.Item(11, 4).Value = 1971.8
$ws.Cells.Item(11, 5).Value = 2001.7
$ws.Cells.Item(11, 6).Value = 586290
$ws.Cells.Item(11, 7).Value = 387913
$ws.Cells.Item(11, 8).Value = 0.5113955964352832
$ws.Cells.Item(11, 9).Value = "GODREJPROP"
$ws.Cells.Item(12, 1).Value = "BHEL"
$ws.Cells.Item(12, 2).Value = 231.2
$ws.Cells.Item(12, 3).Value = 235.6
$ws.Cells.Item(12, 4).Value = 231.01
$ws.Cells.Item(12, 5).Value = 233.5
$ws.Cells.Item(12, 6).Value = 6844860
$ws.Cells.Item(12, 7).Value = 4520518
$ws.Cells.Item(12, 8).Value = 0.5141760302690975
$ws.Cells.Item(12, 9).Value = "BHEL"
$ws.Cells.Item(13, 1).Value = "CUMMINSIND"
$ws.Cells.Item(13, 2).Value = 3947.6
$ws.Cells.Item(13, 3).Value = 4010.3
$ws.Cells.Item(13, 4).Value = 3930
$ws.Cells.Item(13, 5).Value = 3984.8
$ws.Cells.Item(13, 6).Value = 355702
$ws.Cells.Item(13, 7).Value = 245108
$ws.Cells.Item(13, 8).Value = 0.451205183021362
$ws.Cells.Item(13, 9).Value = "CUMMINSIND"
$ws.Cells.Item(14, 1).Value = "JUBLFOOD"
$ws.Cells.Item(14, 2).Value = 610
$ws.Cells.Item(14, 3).Value = 618.65
$ws.Cells.Item(14, 4).Value = 604.15
$ws.Cells.Item(14, 5).Value = 610.2
$ws.Cells.Item(14, 6).Value = 3154893
$ws.Cells.Item(14, 7).Value = 2045929
$ws.Cells.Item(14, 8).Value = 0.5420344498758266
$ws.Cells.Item(14, 9).Value = "JUBLFOOD"
$ws.Cells.Item(15, 1).Value = "ABFRL"
$ws.Cells.Item(15, 2).Value = 85.75
$ws.Cells.Item(15, 3).Value = 85.76
$ws.Cells.Item(15, 4).Value = 83.15
$ws.Cells.Item(15, 5).Value = 85.15
$ws.Cells.Item(15, 6).Value = 9055599
$ws.Cells.Item(15, 7).Value = 5928632
$ws.Cells.Item(15, 8).Value = 0.5274348281357318
$ws.Cells.Item(15, 9).Value = "ABFRL"
$ws.Cells.Item(16, 1).Value = "SRF"
$ws.Cells.Item(16, 2).Value = 2798
$ws.Cells.Item(16, 3).Value = 2850.8
$ws.Cells.Item(16, 4).Value = 2789
$ws.Cells.Item(16, 5).Value = 2814.7
$ws.Cells.Item(16, 6).Value = 221510
$ws.Cells.Item(16, 7).Value = 138769
$ws.Cells.Item(16, 8).Value = 0.5962498828989183
$ws.Cells.Item(16, 9).Value = "SRF"
$ws.Cells.Item(17, 1).Value = "ASTRAL"
$ws.Cells.Item(17, 2).Value = 1373.9
$ws.Cells.Item(17, 3).Value = 1388
$ws.Cells.Item(17, 4).Value = 1364.3
$ws.Cells.Item(17, 5).Value = 1364.9
$ws.Cells.Item(17, 6).Value = 580505
$ws.Cells.Item(17, 7).Value = 376171
$ws.Cells.Item(17, 8).Value = 0.5431944514595756
$ws.Cells.Item(17, 9).Value = "ASTRAL"
$ws.Cells.Item(18, 1).Value = "HFCL"
$ws.Cells.Item(18, 2).Value = 72.46
$ws.Cells.Item(18, 3).Value = 76.19
$ws.Cells.Item(18, 4).Value = 72.3
$ws.Cells.Item(18, 5).Value = 74.2
$ws.Cells.Item(18, 6).Value = 31113853
$ws.Cells.Item(18, 7).Value = 20761048
$ws.Cells.Item(18, 8).Value = 0.4986648554543104
$ws.Cells.Item(18, 9).Value = "HFCL"
$ws.Cells.Item(19, 1).Value = "KFINTECH"
$ws.Cells.Item(19, 2).Value = 1070.8
$ws.Cells.Item(19, 3).Value = 1093.6
$ws.Cells.Item(19, 4).Value = 1056.1
$ws.Cells.Item(19, 5).Value = 1070
$ws.Cells.Item(19, 6).Value = 994466
$ws.Cells.Item(19, 7).Value = 704238
$ws.Cells.Item(19, 8).Value = 0.4121163583902033
$ws.Cells.Item(19, 9).Value = "KFINTECH"

Write-Output "edit complete"
